# Update temp/rh QC check.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row ---
# "data_quality" header is replaced by the new "rh_data_quality" label.
$ws.Range("D1").Value = "rh_data_quality"
$ws.Range("E1").Value = "notes"

# --- Duplicate the existing 24-row QC table (rows 2:25) as a new visit ---
# (row 26 is left blank, new data starts at row 27, mirroring the original layout)
$src = $ws.Range("A2:E25")
$dst = $ws.Range("A27")
$src.Copy($dst)

# The new block is for the 2019-01-09 visit
$ws.Range("B27:B50").Value = 20190109

# --- Update the data_quality / notes values for the new visit ---
$ws.Range("D27").Value = "okay"
$ws.Range("E27").ClearContents()

$ws.Range("D28").Value = "bad"
$ws.Range("E28").Value = "flatlined at 1%"

$ws.Range("D29").Value = "bad"
$ws.Range("E29").Value = "flatlined at 1%"

$ws.Range("D30").Value = "okay"
$ws.Range("E30").Value = "maybe some questionable areas"

$ws.Range("D31").Value = "check"
$ws.Range("E31").Value = "some flatlines at 1%"

$ws.Range("D32").Value = "bad"
$ws.Range("E32").Value = "flatlined at 1%"

$ws.Range("D33").Value = "okay"
$ws.Range("E33").ClearContents()

$ws.Range("D34").Value = "bad"
$ws.Range("E34").Value = "flatlined at 1% prior to 1/2"

$ws.Range("D35").Value = "bad"
$ws.Range("E35").Value = "flatlined at 1%, replaced with logger that was at plot 18"

$ws.Range("D36").Value = "okay"
$ws.Range("E36").ClearContents()

$ws.Range("D37").Value = "okay"
$ws.Range("E37").ClearContents()

$ws.Range("D38").Value = "okay"
$ws.Range("E38").ClearContents()

$ws.Range("D39").Value = "okay"
$ws.Range("E39").ClearContents()

$ws.Range("D40").Value = "okay"
$ws.Range("E40").ClearContents()

$ws.Range("D41").Value = "okay"
$ws.Range("E41").ClearContents()

$ws.Range("D42").Value = "okay"
$ws.Range("E42").ClearContents()

$ws.Range("D43").Value = "okay"
$ws.Range("E43").ClearContents()

# Row 44 (id 18) never had a logger - no C/D values, only a note.
$ws.Range("C44").ClearContents()
$ws.Range("D44").ClearContents()
$ws.Range("E44").Value = "removed during prior visit"

$ws.Range("D45").Value = "okay"
$ws.Range("E45").ClearContents()

$ws.Range("D46").Value = "check"
$ws.Range("E46").Value = "flatlined at 1% prior to 12/18"

$ws.Range("D47").Value = "okay"
$ws.Range("E47").ClearContents()

$ws.Range("D48").Value = "check"
$ws.Range("E48").Value = "flatlined at 1% prior to 12/29"

$ws.Range("D49").Value = "okay"
$ws.Range("E49").ClearContents()

$ws.Range("D50").Value = "bad"
$ws.Range("E50").Value = "flatlined at 1% from 12/18"

# --- Update the sheet view to match where the user ended up ---
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("E38").Select()
